$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D10").Value = 44340
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("P10").Value = 1538
$ws.Range("I11").Value = 'Segunda'
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 18000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 18000
$ws.Range("P11").Value = 1385
$ws.Range("D12").Value = 44333
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 340
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 26000
$ws.Range("M12").Value = 25500
$ws.Range("P12").Value = 1962
$ws.Range("I13").Value = 'Segunda'
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 23000
$ws.Range("L13").Value = 23000
$ws.Range("M13").Value = 23000
$ws.Range("P13").Value = 1769
$ws.Range("D14").Value = 44445
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 790
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13494
$ws.Range("P14").Value = 1038
$ws.Range("I15").Value = 'Segunda'
$ws.Range("J15").Value = 340
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 885
$ws.Range("D16").Value = 44648
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 610
$ws.Range("K16").Value = 16000
$ws.Range("L16").Value = 17000
$ws.Range("M16").Value = 16500
$ws.Range("P16").Value = 1269
$ws.Range("D17").Value = 44767
$ws.Range("J17").Value = 700
$ws.Range("K17").Value = 17000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = 17500
$ws.Range("P17").Value = 1346
$ws.Range("D18").Value = 44627
$ws.Range("J18").Value = 790
$ws.Range("K18").Value = 14000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 14494
$ws.Range("P18").Value = 1115
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 340
$ws.Range("K19").Value = 13000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 13000
$ws.Range("P19").Value = 1000
$ws.Range("D20").Value = 44914
$ws.Range("I20").Value = 'Primera'
$ws.Range("J20").Value = 520
$ws.Range("L20").Value = 14000
$ws.Range("M20").Value = 13500
$ws.Range("P20").Value = 1038
$ws.Range("D21").Value = 44172
$ws.Range("J21").Value = 430
$ws.Range("K21").Value = 30000
$ws.Range("L21").Value = 30000
$ws.Range("M21").Value = 30000
$ws.Range("P21").Value = 2308
$ws.Range("D22").Value = 44263
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 40000
$ws.Range("L22").Value = 40000
$ws.Range("M22").Value = 40000
$ws.Range("N22").Value = '$/caja 15 kilos'
$ws.Range("P22").Value = 2667
$ws.Range("Q22").Value = 15
$ws.Range("D23").Value = 44200
$ws.Range("J23").Value = 520
$ws.Range("K23").Value = 30000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 30000
$ws.Range("N23").Value = '$/caja 13 kilos'
$ws.Range("P23").Value = 2308
$ws.Range("Q23").Value = 13
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 340
$ws.Range("K24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("M24").Value = 25000
$ws.Range("P24").Value = 1923
$ws.Range("D25").Value = 44417
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 790
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 1115
$ws.Range("I26").Value = 'Segunda'
$ws.Range("J26").Value = 340
$ws.Range("K26").Value = 13000
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = 13000
$ws.Range("P26").Value = 1000
$ws.Range("D27").Value = 44690
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 790
$ws.Range("K27").Value = 12000
$ws.Range("M27").Value = 12494
$ws.Range("P27").Value = 961
$ws.Range("D28").Value = 44389
$ws.Range("J28").Value = 700
$ws.Range("K28").Value = 19000
$ws.Range("L28").Value = 20000
$ws.Range("M28").Value = 19500
$ws.Range("P28").Value = 1500
$ws.Range("I29").Value = 'Segunda'
$ws.Range("J29").Value = 340
$ws.Range("K29").Value = 17000
$ws.Range("L29").Value = 17000
$ws.Range("M29").Value = 17000
$ws.Range("P29").Value = 1308
$ws.Range("D30").Value = 44410
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 790
$ws.Range("K30").Value = 15000
$ws.Range("L30").Value = 16000
$ws.Range("M30").Value = 15494
$ws.Range("P30").Value = 1192
$ws.Range("I31").Value = 'Segunda'
$ws.Range("J31").Value = 340
$ws.Range("K31").Value = 13000
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 13000
$ws.Range("P31").Value = 1000
$ws.Range("D32").Value = 44522
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 790
$ws.Range("K32").Value = 16000
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 16987
$ws.Range("P32").Value = 1307
$ws.Range("I33").Value = 'Segunda'
$ws.Range("J33").Value = 360
$ws.Range("K33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("M33").Value = 15000
$ws.Range("P33").Value = 1154
$ws.Range("D34").Value = 44809
$ws.Range("I34").Value = 'Primera'
$ws.Range("J34").Value = 350
$ws.Range("K34").Value = 13000
$ws.Range("M34").Value = 14143
$ws.Range("P34").Value = 1088
$ws.Range("I35").Value = 'Segunda'
$ws.Range("J35").Value = 160
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 10000
$ws.Range("P35").Value = 769
$ws.Range("D36").Value = 44949
$ws.Range("I36").Value = 'Primera'
$ws.Range("J36").Value = 340
$ws.Range("K36").Value = 24000
$ws.Range("L36").Value = 25000
$ws.Range("M36").Value = 24500
$ws.Range("P36").Value = 1885
$ws.Range("D37").Value = 44998
$ws.Range("J37").Value = 790
$ws.Range("K37").Value = 17000
$ws.Range("L37").Value = 18000
$ws.Range("M37").Value = 17506
$ws.Range("P37").Value = 1347
